$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected range to Text format so the numeric-looking / percent-looking
# strings are stored as literal text (matching the workbook's existing inlineStr
# text cells) instead of being auto-converted to numbers/percentages by Excel's
# automatic data recognition.
$dataRange = $ws.Range("D2:E50")
$dataRange.NumberFormat = "@"

$ws.Range("E2").Value = "0.76%"
$ws.Range("D3").Value = "26.86"
$ws.Range("E3").Value = "-1.78%"
$ws.Range("D4").Value = "4.684"
$ws.Range("E4").Value = "3.39%"
$ws.Range("D5").Value = "0.05989"
$ws.Range("E5").Value = "1.66%"
$ws.Range("D6").Value = "6.660"
$ws.Range("E6").Value = "0.51%"
$ws.Range("D7").Value = "0.8555"
$ws.Range("E7").Value = "-0.22%"
$ws.Range("D8").Value = "0.9217"
$ws.Range("E8").Value = "-0.91%"
$ws.Range("D9").Value = "0.1388"
$ws.Range("D10").Value = "0.05147"
$ws.Range("E10").Value = "40.43%"
$ws.Range("D11").Value = "0.07009"
$ws.Range("E11").Value = "-1.17%"
$ws.Range("D12").Value = "0.03047"
$ws.Range("E12").Value = "-5.70%"
$ws.Range("D13").Value = "0.09130"
$ws.Range("E13").Value = "-0.82%"
$ws.Range("D14").Value = "0.001526"
$ws.Range("E14").Value = "-0.84%"
$ws.Range("D15").Value = "0.0006041"
$ws.Range("E15").Value = "-0.45%"
$ws.Range("D16").Value = "0.006107"
$ws.Range("E16").Value = "1.31%"
$ws.Range("D17").Value = "3.453"
$ws.Range("E17").Value = "-1.76%"
$ws.Range("E18").Value = "-1.66%"
$ws.Range("E19").Value = "-1.52%"
$ws.Range("D20").Value = "0.3110"
$ws.Range("E20").Value = "1.69%"
$ws.Range("D22").Value = "4.135"
$ws.Range("E22").Value = "7.51%"
$ws.Range("D23").Value = "0.04232"
$ws.Range("E23").Value = "0.43%"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").Value = "-0.60%"
$ws.Range("D25").Value = "0.004036"
$ws.Range("E25").Value = "-5.82%"
$ws.Range("E26").Value = "-0.07%"
$ws.Range("E27").Value = "13.36%"
$ws.Range("D40").Value = "0.03833"
$ws.Range("E40").Value = "-0.04%"
$ws.Range("E41").Value = "1.29%"
$ws.Range("D42").Value = "0.003806"
$ws.Range("E42").Value = "-38.49%"
$ws.Range("D43").Value = "0.002417"
$ws.Range("E43").Value = "9.93%"
$ws.Range("D44").Value = "0.01504"
$ws.Range("E44").Value = "32.36%"
$ws.Range("D45").Value = "0.00005116"
$ws.Range("E45").Value = "-6.40%"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("E47").Value = "-30.60%"
$ws.Range("D48").Value = "0.1503"
$ws.Range("E48").Value = "35.93%"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.04%"


# Restore the original (default/"Normal") style on the range so we don't leave
# a stray NumberFormat/style behind on cells that originally had no explicit
# style applied.
$dataRange.Style = "Normal"

Write-Output "Updated $($dataRange.Address()) crypto price/volume figures."
